$wb = $excel.ActiveWorkbook

# Update the "all_tools" sheet: column C (warnings count) rows 10-12 change from 39 to 40
$wsAllTools = $wb.Worksheets.Item("all_tools")
$wsAllTools.Range("C10").Value = 40
$wsAllTools.Range("C11").Value = 40
$wsAllTools.Range("C12").Value = 40

# Update the "infer" sheet: recalculated correlation statistics for rows 10-12 (F:I)
$wsInfer = $wb.Worksheets.Item("infer")

$wsInfer.Range("F10").Value = -0.2465459984594313
$wsInfer.Range("G10").Value = 0.03317524973370743
$wsInfer.Range("H10").Value = -0.3098303775675377
$wsInfer.Range("I10").Value = 0.02855238451552488

$wsInfer.Range("F11").Value = -0.02257043088830945
$wsInfer.Range("G11").Value = 0.8463887047098356
$wsInfer.Range("H11").Value = -0.02588543710805493
$wsInfer.Range("I11").Value = 0.8583790451496913

$wsInfer.Range("F12").Value = -0.007967906165899995
$wsInfer.Range("G12").Value = 0.9479195406770614
$wsInfer.Range("H12").Value = -0.007746944922644224
$wsInfer.Range("I12").Value = 0.9574176555277638
